$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B..AC (id column A is left untouched per row)
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC")

function Swap-Rows($rowA, $rowB) {
    foreach ($col in $cols) {
        $cellA = $ws.Range($col + $rowA)
        $cellB = $ws.Range($col + $rowB)
        $valA = $cellA.Value()
        $valB = $cellB.Value()
        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}

# Swap the mismatched row pairs back into the correct order
Swap-Rows 9 10
Swap-Rows 36 37
Swap-Rows 76 77
Swap-Rows 99 100

# Append the new match result as row 157
$newRow = 157
$ws.Cells.Item(2, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)
$ws.Cells.Item(2, 5).Copy()
$ws.Cells.Item($newRow, 5).PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value = 155
$ws.Cells.Item($newRow, 2).Value = 7952747
$ws.Cells.Item($newRow, 3).Value = "Bosnia Herzegovina Premier Liga"
$ws.Cells.Item($newRow, 4).Value = "Bosnia  Herzegovina Premier Liga"
$ws.Cells.Item($newRow, 5).Value = 45396.625
$ws.Cells.Item($newRow, 6).Value = "Zeljeznicar"
$ws.Cells.Item($newRow, 7).Value = "FK Sarajevo"
$ws.Cells.Item($newRow, 8).Value = 0
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = "D"
$ws.Cells.Item($newRow, 11).Value = 2.4
$ws.Cells.Item($newRow, 12).Value = 3.2
$ws.Cells.Item($newRow, 13).Value = 2.6
$ws.Cells.Item($newRow, 14).Value = 2.4
$ws.Cells.Item($newRow, 15).Value = 3.2
$ws.Cells.Item($newRow, 16).Value = 2.625
$ws.Cells.Item($newRow, 17).Value = 0
$ws.Cells.Item($newRow, 18).Value = 1.8
$ws.Cells.Item($newRow, 19).Value = 2
$ws.Cells.Item($newRow, 20).Value = 2.25
$ws.Cells.Item($newRow, 21).Value = 1.95
$ws.Cells.Item($newRow, 22).Value = 1.85
$ws.Cells.Item($newRow, 23).Value = -1
$ws.Cells.Item($newRow, 24).Value = 2.2
$ws.Cells.Item($newRow, 25).Value = -1
$ws.Cells.Item($newRow, 26).Value = 0
$ws.Cells.Item($newRow, 27).Value = -0
$ws.Cells.Item($newRow, 28).Value = -1
$ws.Cells.Item($newRow, 29).Value = 0.8500000000000001
